{"js": "// Update the date heading paragraph (first paragraph of the body, outside the table).\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paras.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\nif (dateParagraph.text.trim() === \"2024-10-21 Monday\") {\n  dateParagraph.insertText(\"2024-10-22 Tuesday\", \"Replace\");\n}\n\n// Update the division-problem answers held in the first table of the document.\n// Values are addressed positionally (row, col) rather than by searching for the\n// old text, because some old/new values collide across cells (e.g. \"99\u00f77=14, 1\"\n// is both a target of one cell and the prior value of another).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map of {row, col, newValue} \u2014 rows are the 5 data rows of the 20-row table\n// (every 4th row starting at 0), columns are 0-4 left to right.\nconst updates = [\n  { row: 0, col: 0, value: \"96\u00f77=13, 5\" },\n  { row: 0, col: 1, value: \"99\u00f77=14, 1\" },\n  { row: 0, col: 2, value: \"82\u00f72=41, 0\" },\n  { row: 0, col: 3, value: \"95\u00f74=23, 3\" },\n  { row: 0, col: 4, value: \"24\u00f79=2, 6\" },\n\n  { row: 4, col: 0, value: \"77\u00f79=8, 5\" },\n  { row: 4, col: 1, value: \"86\u00f79=9, 5\" },\n  { row: 4, col: 2, value: \"42\u00f72=21, 0\" },\n  { row: 4, col: 3, value: \"90\u00f77=12, 6\" },\n  { row: 4, col: 4, value: \"14\u00f78=1, 6\" },\n\n  { row: 8, col: 0, value: \"18\u00f76=3, 0\" },\n  { row: 8, col: 1, value: \"88\u00f72=44, 0\" },\n  { row: 8, col: 2, value: \"58\u00f72=29, 0\" },\n  { row: 8, col: 3, value: \"37\u00f78=4, 5\" },\n  { row: 8, col: 4, value: \"78\u00f76=13, 0\" },\n\n  { row: 12, col: 0, value: \"72\u00f79=8, 0\" },\n  { row: 12, col: 1, value: \"92\u00f74=23, 0\" },\n  { row: 12, col: 2, value: \"44\u00f79=4, 8\" },\n  { row: 12, col: 3, value: \"73\u00f72=36, 1\" },\n  { row: 12, col: 4, value: \"95\u00f79=10, 5\" },\n\n  { row: 16, col: 0, value: \"24\u00f75=4, 4\" },\n  { row: 16, col: 1, value: \"14\u00f76=2, 2\" },\n  { row: 16, col: 2, value: \"56\u00f74=14, 0\" },\n  { row: 16, col: 3, value: \"75\u00f74=18, 3\" },\n  { row: 16, col: 4, value: \"70\u00f75=14, 0\" },\n];\n\nfor (const u of updates) {\n  const cell = table.getCell(u.row, u.col);\n  cell.value = u.value;\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading (the document's first paragraph, outside the table).\n$d = $word.ActiveDocument\n\n$dateParagraph = $d.Paragraphs.Item(1)\nif ($dateParagraph.Range.Text.TrimEnd([char]13, [char]7) -eq \"2024-10-21 Monday\") {\n    $dateParagraph.Range.Text = \"2024-10-22 Tuesday\"\n}\n\n# Update the division-problem answers held in the first table of the document.\n# Cells are addressed positionally (1-based row, col) rather than via text\n# search-and-replace, because some old/new values collide across cells\n# (e.g. \"99\u00f77=14, 1\" is both the new value of one cell and the prior value\n# of another), which would make a naive global Find/Replace order-dependent.\n$t = $d.Tables.Item(1)\n\n# Data lives in rows 1, 5, 9, 13, 17 (1-based; every 4th row) of the 20-row\n# table; columns 1-5 left to right.\n$t.Cell(1, 1).Range.Text  = \"96\u00f77=13, 5\"\n$t.Cell(1, 2).Range.Text  = \"99\u00f77=14, 1\"\n$t.Cell(1, 3).Range.Text  = \"82\u00f72=41, 0\"\n$t.Cell(1, 4).Range.Text  = \"95\u00f74=23, 3\"\n$t.Cell(1, 5).Range.Text  = \"24\u00f79=2, 6\"\n\n$t.Cell(5, 1).Range.Text  = \"77\u00f79=8, 5\"\n$t.Cell(5, 2).Range.Text  = \"86\u00f79=9, 5\"\n$t.Cell(5, 3).Range.Text  = \"42\u00f72=21, 0\"\n$t.Cell(5, 4).Range.Text  = \"90\u00f77=12, 6\"\n$t.Cell(5, 5).Range.Text  = \"14\u00f78=1, 6\"\n\n$t.Cell(9, 1).Range.Text  = \"18\u00f76=3, 0\"\n$t.Cell(9, 2).Range.Text  = \"88\u00f72=44, 0\"\n$t.Cell(9, 3).Range.Text  = \"58\u00f72=29, 0\"\n$t.Cell(9, 4).Range.Text  = \"37\u00f78=4, 5\"\n$t.Cell(9, 5).Range.Text  = \"78\u00f76=13, 0\"\n\n$t.Cell(13, 1).Range.Text = \"72\u00f79=8, 0\"\n$t.Cell(13, 2).Range.Text = \"92\u00f74=23, 0\"\n$t.Cell(13, 3).Range.Text = \"44\u00f79=4, 8\"\n$t.Cell(13, 4).Range.Text = \"73\u00f72=36, 1\"\n$t.Cell(13, 5).Range.Text = \"95\u00f79=10, 5\"\n\n$t.Cell(17, 1).Range.Text = \"24\u00f75=4, 4\"\n$t.Cell(17, 2).Range.Text = \"14\u00f76=2, 2\"\n$t.Cell(17, 3).Range.Text = \"56\u00f74=14, 0\"\n$t.Cell(17, 4).Range.Text = \"75\u00f74=18, 3\"\n$t.Cell(17, 5).Range.Text = \"70\u00f75=14, 0\"\n"}
